$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New scan record appended as row 2 of the QR Scanner log.
# A2 ("228723") and C2 ("05/17/2025") look like a number / date to Excel's
# auto-detection, so they're entered with a leading apostrophe to force
# text, then the cell style is reset to "Normal" so no extra
# number-format/quote-prefix styling is left behind on the cell.
$ws.Range("A2").Value = "'228723"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "Biochemistry"

$ws.Range("C2").Value = "'05/17/2025"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").Value = "22:44:27"
$ws.Range("E2").Value = "Manual"
$ws.Range("F2").Value = "231249@med.asu.edu.eg"
